$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 14 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 14
}

# Update the "Notified Production (MW)" values for rows 22-89 (column B)
$bValues = @{
    22 = 1;    23 = 1;    24 = 2;    25 = 2;    26 = 5;    27 = 6;
    28 = 9;    29 = 18;   30 = 119;  31 = 158;  32 = 216;  33 = 262;
    34 = 537;  35 = 613;  36 = 694;  37 = 779;  38 = 1014; 39 = 1097;
    40 = 1193; 41 = 1279; 42 = 1467; 43 = 1538; 44 = 1593; 45 = 1643;
    46 = 1701; 47 = 1727; 48 = 1746; 49 = 1751; 50 = 1741; 51 = 1733;
    52 = 1715; 53 = 1695; 54 = 1638; 55 = 1606; 56 = 1563; 57 = 1521;
    58 = 1386; 59 = 1325; 60 = 1254; 61 = 1192; 62 = 985;  63 = 902;
    64 = 815;  65 = 730;  66 = 511;  67 = 408;  68 = 315;  69 = 241;
    70 = 93;   71 = 48;   72 = 34;   73 = 27;   74 = 13;   75 = 11;
    76 = 6;    77 = 5;    78 = 5;    79 = 5;    80 = 5;    81 = 3;
    82 = 1;    83 = 1;    84 = 1;    85 = 1;    86 = 1;    87 = 1;
    88 = 1;    89 = 1
}

foreach ($r in $bValues.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $bValues[$r]
}
